$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Capture original values before overwriting (3-cycle rotation among rows 2, 3, 5)
$d2 = $ws.Range("D2").Value2
$m2 = $ws.Range("M2").Value2
$d3 = $ws.Range("D3").Value2
$m3 = $ws.Range("M3").Value2
$d5 = $ws.Range("D5").Value2
$m5 = $ws.Range("M5").Value2

# Apply rotation: row2 <- old row5, row3 <- old row2, row5 <- old row3
$ws.Range("D2").Value2 = $d5
$ws.Range("M2").Value2 = $m5

$ws.Range("D3").Value2 = $d2
$ws.Range("M3").Value2 = $m2

$ws.Range("D5").Value2 = $d3
$ws.Range("M5").Value2 = $m3
